# Append real-time data rows through June 2023 (commit: "updated real time data to 2023 June")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (date number format, etc.) of the last existing row (664)
# down into the new rows (665:678) so the new date cells keep the same style.
$ws.Range("A664:C664").Copy() | Out-Null
$ws.Range("A665:C678").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(665, 1).Value = 44682
$ws.Cells.Item(665, 2).Value = 291.474
$ws.Cells.Item(665, 3).Value = 292.289

$ws.Cells.Item(666, 1).Value = 44713
$ws.Cells.Item(666, 2).Value = 295.328
$ws.Cells.Item(666, 3).Value = 294.354

$ws.Cells.Item(667, 1).Value = 44743
$ws.Cells.Item(667, 2).Value = 295.271
$ws.Cells.Item(667, 3).Value = 295.275

$ws.Cells.Item(668, 1).Value = 44774
$ws.Cells.Item(668, 2).Value = 295.62
$ws.Cells.Item(668, 3).Value = 296.95

$ws.Cells.Item(669, 1).Value = 44805
$ws.Cells.Item(669, 2).Value = 296.761
$ws.Cells.Item(669, 3).Value = 298.66

$ws.Cells.Item(670, 1).Value = 44835
$ws.Cells.Item(670, 2).Value = 298.062
$ws.Cells.Item(670, 3).Value = 299.471

$ws.Cells.Item(671, 1).Value = 44866
$ws.Cells.Item(671, 2).Value = 298.349
$ws.Cells.Item(671, 3).Value = 300.066

$ws.Cells.Item(672, 1).Value = 44896
$ws.Cells.Item(672, 2).Value = 298.112
$ws.Cells.Item(672, 3).Value = 300.974

$ws.Cells.Item(673, 1).Value = 44927
$ws.Cells.Item(673, 2).Value = 300.536
$ws.Cells.Item(673, 3).Value = 302.702

$ws.Cells.Item(674, 1).Value = 44958
$ws.Cells.Item(674, 2).Value = 301.648
$ws.Cells.Item(674, 3).Value = 304.07

$ws.Cells.Item(675, 1).Value = 44986
$ws.Cells.Item(675, 2).Value = 301.808
$ws.Cells.Item(675, 3).Value = 305.24

$ws.Cells.Item(676, 1).Value = 45017
$ws.Cells.Item(676, 2).Value = 302.918
$ws.Cells.Item(676, 3).Value = 306.489

$ws.Cells.Item(677, 1).Value = 45047
$ws.Cells.Item(677, 2).Value = 303.294
$ws.Cells.Item(677, 3).Value = 307.824

$ws.Cells.Item(678, 1).Value = 45078
$ws.Cells.Item(678, 2).Value = 303.841
$ws.Cells.Item(678, 3).Value = 308.309
